$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated / reordered data (rows 2-16), sorted descending by total_registros,
# with a new row appended at the end (PIÑIN MACHUCA GIANCARLO, 1)
$names = @(
    "ALBIRENA GARCIA ANGEELO ALONSO",
    "URRIOLA ARISMENDIZ INGRID MARYURI",
    "MANUEL LEUNARDO PRADO BAILON",
    "CORDOVA CARMEN ANGIE NATALLY",
    "MARYURI OJEDA VALLE",
    "ATOCHE PALACIOS LUIS ANGEL",
    "ALAMA NIMA CLARITZA MABEL",
    "AGURTO ORDINOLA LISBET JAQUELIN",
    "RUIDIAS FRIAS MELISSA VICTORIA",
    "VEGA ROBLEDO FERNANDO ERNESTO",
    "JUAREZ CARMEN PIERRE ALEXANDER",
    "ROMAN GALECIO MARITZA DEL PILAR",
    "BERNAOLA CARMEN ZUMIKO YASHURY",
    "CARREÑO PALACIOS KATHERINE DE LOS MILAGROS",
    "PIÑIN MACHUCA GIANCARLO"
)

$totals = @(208, 197, 193, 184, 184, 182, 179, 174, 169, 159, 157, 150, 143, 120, 1)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $totals[$i]
}
